# Updated cryptos list on Sun Oct  8 10:50:37 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for each
# coin row, and reorders the Aave/RenderToken and MXToken/RocketPoolETH
# pairs (rows 44-47) to reflect the newly scraped ranking snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds values that look numeric (e.g. "0.520", "4.00",
# "27.312.28") but must stay as literal text, matching the original
# inline-string cells. Force the column to Text format first so Excel does
# not "helpfully" convert them to numbers and drop significant trailing
# zeros / thousands separators.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "27.312.28"
$ws.Range("E2").Value = "  -2.30%  "
# Row 3
$ws.Range("D3").Value = "1.615.83"
$ws.Range("E3").Value = "  -1.50%  "
# Row 4
$ws.Range("D4").Value = "0.993"
$ws.Range("E4").Value = "  -0.75%  "
# Row 5
$ws.Range("D5").Value = "208.88"
$ws.Range("E5").Value = "  -1.94%  "
# Row 6
$ws.Range("D6").Value = "0.520"
$ws.Range("E6").Value = "  -0.76%  "
# Row 7
$ws.Range("D7").Value = "0.994"
$ws.Range("E7").Value = "  -0.60%  "
# Row 8
$ws.Range("D8").Value = "23.17"
$ws.Range("E8").Value = "  -1.55%  "
# Row 9
$ws.Range("D9").Value = "0.255"
$ws.Range("E9").Value = "  -2.33%  "
# Row 10
$ws.Range("D10").Value = "0.0610"
$ws.Range("E10").Value = "  -0.85%  "
# Row 11
$ws.Range("D11").Value = "0.0875"
$ws.Range("E11").Value = "  -0.61%  "
# Row 12
$ws.Range("D12").Value = "1.854.88"
$ws.Range("E12").Value = "  -0.96%  "
# Row 13
$ws.Range("D13").Value = "1.619.20"
$ws.Range("E13").Value = "  -1.19%  "
# Row 14
$ws.Range("D14").Value = "4.00"
$ws.Range("E14").Value = "  -2.31%  "
# Row 15
$ws.Range("D15").Value = "0.557"
$ws.Range("E15").Value = "  -2.82%  "
# Row 16
$ws.Range("D16").Value = "64.99"
$ws.Range("E16").Value = "  -1.41%  "
# Row 17
$ws.Range("D17").Value = "27.578.46"
$ws.Range("E17").Value = "  -1.35%  "
# Row 18
$ws.Range("D18").Value = "227.59"
$ws.Range("E18").Value = "  -2.50%  "
# Row 19
$ws.Range("D19").Value = "0.0₃0718"
$ws.Range("E19").Value = "  -0.90%  "
# Row 20
$ws.Range("D20").Value = "7.56"
$ws.Range("E20").Value = "  -0.81%  "
# Row 21
$ws.Range("D21").Value = "0.997"
$ws.Range("E21").Value = "  -0.38%  "
# Row 22
$ws.Range("D22").Value = "4.28"
$ws.Range("E22").Value = "  -1.92%  "
# Row 23
$ws.Range("D23").Value = "10.02"
$ws.Range("E23").Value = "  -6.95%  "
# Row 24
$ws.Range("D24").Value = "2.04"
$ws.Range("E24").Value = "  -2.01%  "
# Row 25
$ws.Range("D25").Value = "153.78"
$ws.Range("E25").Value = "  +1.71%  "
# Row 26
$ws.Range("D26").Value = "6.87"
$ws.Range("E26").Value = "  -1.57%  "
# Row 27
$ws.Range("E27").Value = "  -0.50%  "
# Row 28
$ws.Range("D28").Value = "15.42"
$ws.Range("E28").Value = "  -1.67%  "
# Row 29
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.36%  "
# Row 30
$ws.Range("E30").Value = "  -1.28%  "
# Row 31
$ws.Range("D31").Value = "0.0479"
$ws.Range("E31").Value = "  -0.99%  "
# Row 32
$ws.Range("D32").Value = "3.39"
$ws.Range("E32").Value = "  +1.27%  "
# Row 33
$ws.Range("D33").Value = "3.06"
$ws.Range("E33").Value = "  -2.04%  "
# Row 34
$ws.Range("D34").Value = "1.381.73"
$ws.Range("E34").Value = "  -2.20%  "
# Row 35
$ws.Range("D35").Value = "1.56"
$ws.Range("E35").Value = "  -1.07%  "
# Row 36
$ws.Range("D36").Value = "0.996"
$ws.Range("E36").Value = "  +10.03%  "
# Row 37
$ws.Range("E37").Value = "  -1.19%  "
# Row 38
$ws.Range("E38").Value = "  +0.09%  "
# Row 39
$ws.Range("D39").Value = "0.552"
$ws.Range("E39").Value = "  -1.19%  "
# Row 40
$ws.Range("D40").Value = "0.848"
$ws.Range("E40").Value = "  -4.06%  "
# Row 41
$ws.Range("D41").Value = "1.02"
$ws.Range("E41").Value = "  -0.54%  "
# Row 42
$ws.Range("D42").Value = "0.996"
$ws.Range("E42").Value = "  -0.48%  "
# Row 43
$ws.Range("D43").Value = "5.45"
$ws.Range("E43").Value = "  -1.29%  "

# Rows 44-45 swap: RenderToken <-> Aave (with freshly scraped price/volume)
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "65.29"
$ws.Range("E44").Value = "  -1.84%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "1.81"
$ws.Range("E45").Value = "  -3.47%  "

# Rows 46-47 swap: RocketPoolETH <-> MXToken (with freshly scraped price/volume)
$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D46").Value = "2.17"
$ws.Range("E46").Value = "  -1.55%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.765.06"
$ws.Range("E47").Value = "  -0.92%  "

# Row 48
$ws.Range("D48").Value = "87.55"
$ws.Range("E48").Value = "  -0.31%  "
# Row 49
$ws.Range("D49").Value = "0.0₆0103"
$ws.Range("E49").Value = "  -2.24%  "
# Row 50
$ws.Range("E50").Value = "  +0.40%  "
# Row 51
$ws.Range("D51").Value = "0.0502"
$ws.Range("E51").Value = "  -0.88%  "
